$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'59.025.49"
$ws.Range("E2").Value = "  +1.99%  "
$ws.Range("D3").Value = "'2.582.80"
$ws.Range("E3").Value = "  +0.64%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "'520.92"
$ws.Range("E5").Value = "  +0.99%  "
$ws.Range("D6").Value = "'138.85"
$ws.Range("E6").Value = "  -2.48%  "
$ws.Range("D8").Value = "'0.563"
$ws.Range("E8").Value = "  -0.40%  "
$ws.Range("D9").Value = "'2.593.65"
$ws.Range("E9").Value = "  +0.49%  "
$ws.Range("D10").Value = "'6.56"
$ws.Range("E10").Value = "  -0.20%  "
$ws.Range("E11").Value = "  +0.03%  "
$ws.Range("E12").Value = "  +1.73%  "
$ws.Range("E13").Value = "  +3.25%  "
$ws.Range("D14").Value = "'3.041.26"
$ws.Range("E14").Value = "  +0.76%  "
$ws.Range("D15").Value = "'58.816.71"
$ws.Range("E15").Value = "  +1.62%  "
$ws.Range("D16").Value = "'20.37"
$ws.Range("E16").Value = "  +0.67%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "'2.610.29"
$ws.Range("E17").Value = "  +1.97%  "
$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").Value = "'0.0000133"
$ws.Range("E18").Value = "  -0.47%  "
$ws.Range("D19").Value = "'336.75"
$ws.Range("E19").Value = "  -0.70%  "
$ws.Range("D20").Value = "'4.28"
$ws.Range("E20").Value = "  +0.06%  "
$ws.Range("D21").Value = "'10.08"
$ws.Range("E21").Value = "  -1.06%  "
$ws.Range("D22").Value = "'6.51"
$ws.Range("E22").Value = "  +3.23%  "
$ws.Range("E23").Value = "  +0.04%  "
$ws.Range("D24").Value = "'65.94"
$ws.Range("E24").Value = "  +0.88%  "
$ws.Range("E25").Value = "  +0.93%  "
$ws.Range("D26").Value = "'0.402"
$ws.Range("E26").Value = "  +0.48%  "
$ws.Range("D27").Value = "'0.999"
$ws.Range("E27").Value = "  +0.20%  "
$ws.Range("E28").Value = "  +0.63%  "
$ws.Range("E29").Value = "  +0.07%  "
$ws.Range("D30").Value = "'0.0₃0722"
$ws.Range("E30").Value = "  -2.54%  "
$ws.Range("E31").Value = "  -5.24%  "
$ws.Range("E32").Value = "  +0.17%  "
$ws.Range("D33").Value = "'18.66"
$ws.Range("E33").Value = "  +0.16%  "
$ws.Range("D34").Value = "'149.15"
$ws.Range("E34").Value = "  -0.43%  "
$ws.Range("E35").Value = "  -0.12%  "
$ws.Range("E36").Value = "  -2.17%  "
$ws.Range("D37").Value = "'36.78"
$ws.Range("E38").Value = "  +0.91%  "
$ws.Range("D39").Value = "'0.822"
$ws.Range("E39").Value = "  -0.87%  "
$ws.Range("D40").Value = "'0.806"
$ws.Range("E40").Value = "  -7.50%  "
$ws.Range("E41").Value = "  -0.42%  "
$ws.Range("E42").Value = "  +0.16%  "
$ws.Range("D43").Value = "'270.49"
$ws.Range("E43").Value = "  +0.32%  "
$ws.Range("D44").Value = "'10.76"
$ws.Range("E44").Value = "  +0.85%  "
$ws.Range("E45").Value = "  +0.21%  "
$ws.Range("D46").Value = "'0.588"
$ws.Range("E46").Value = "  +0.41%  "
$ws.Range("D47").Value = "'0.0516"
$ws.Range("E47").Value = "  -0.76%  "
$ws.Range("D48").Value = "'18.39"
$ws.Range("E48").Value = "  -1.62%  "
$ws.Range("D49").Value = "'1.963.76"
$ws.Range("E49").Value = "  -0.81%  "
$ws.Range("E50").Value = "  -1.74%  "
$ws.Range("E51").Value = "  -0.21%  "
